# Add Screener "Runs" history tab (S26_G02) rows to the sprint tasks sheet:
# - S26_G02_TB001 (backend run history APIs)
# - S26_G02_TF001 (frontend Runs tab)
# - S26_G02_TT001 (tests)
# - S26_G02_TB002 (UTC timestamp normalization)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-formatted (non-wrapped / default aligned) text value
# into a cell, matching the "default style" cells already used elsewhere in
# this sheet for freshly-authored rows (as opposed to the wrapped/top-aligned
# style inherited from the column defaults).
function Set-PlainCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.HorizontalAlignment = -4132   # xlGeneral
    $c.VerticalAlignment = -4107     # xlBottom
    $c.WrapText = $false
    $c.Value = $text
}

# Row 255: S26_G02_TB001
Set-PlainCell 255 1 "S26"
Set-PlainCell 255 2 "G02"
Set-PlainCell 255 3 "Screener: run history & retention"
Set-PlainCell 255 4 "S26_G02_TB001"
Set-PlainCell 255 5 "Backend: Add Screener V3 run history APIs (list/delete/cleanup) and enrich run read model with query metadata (targets, variables, DSL)."
Set-PlainCell 255 6 "Retention cleanup applies max_days then max_runs; heavy rows remain opt-in via include_rows."
Set-PlainCell 255 7 "implemented"
Set-PlainCell 255 8 "Added /api/screener-v3/runs list/delete + /runs/cleanup; run payload now includes include_holdings/group_ids/variables/condition_dsl for reloading past runs."

# Row 256: S26_G02_TF001
Set-PlainCell 256 1 "S26"
Set-PlainCell 256 2 "G02"
Set-PlainCell 256 3 "Screener: run history & retention"
Set-PlainCell 256 4 "S26_G02_TF001"
Set-PlainCell 256 5 "Frontend: Add Results/Runs tabs on Screener; runs table with View/Load/Delete; retention controls and optional auto-cleanup."
Set-PlainCell 256 6 "Retention settings stored in localStorage (user-managed)."
Set-PlainCell 256 7 "implemented"
Set-PlainCell 256 8 "Screener right panel now supports Runs tab; clicking Run # opens history; retention cleanup can be run manually or automatically."

# Row 257: S26_G02_TT001 (no deviations/F value)
Set-PlainCell 257 1 "S26"
Set-PlainCell 257 2 "G02"
Set-PlainCell 257 3 "Screener: run history & retention"
Set-PlainCell 257 4 "S26_G02_TT001"
Set-PlainCell 257 5 "Tests: Add coverage for Screener run list/delete/cleanup endpoints and ensure lint/build passes."
Set-PlainCell 257 7 "implemented"
Set-PlainCell 257 8 "Added backend/tests/test_screener_runs_api.py and verified ruff/pytest + frontend build."

# Row 258: S26_G02_TB002 (no deviations/F value)
Set-PlainCell 258 1 "S26"
Set-PlainCell 258 2 "G02"
Set-PlainCell 258 3 "Screener: run history & retention"
Set-PlainCell 258 4 "S26_G02_TB002"
Set-PlainCell 258 5 "Backend: Ensure screener run timestamps are UTC-aware so UI displays correct local time (IST +5:30)."
Set-PlainCell 258 7 "implemented"
Set-PlainCell 258 8 "Normalized ScreenerRun created_at/started_at/finished_at to UTC when serializing."
